# Auto-generated update of market-price-derived columns (H:N) for the Leve profit
# tables on each job sheet, reflecting the latest scheduled price-fetch run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 1776.7307
$ws.Range("I113").Value = 1439.8
$ws.Range("J113").Value = 1856.9524
$ws.Range("K113").Value = 1439.8
$ws.Range("L113").Value = 1856.9524
$ws.Range("M113").Value = 1814.2
$ws.Range("N113").Value = -8364.9524
$ws.Range("H129").Value = 716.83
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 716.83
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 2150.49
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -12150.49
$ws.Range("H132").Value = 4193.1304
$ws.Range("I132").Value = 4124.6
$ws.Range("J132").Value = 4650
$ws.Range("K132").Value = 12373.8
$ws.Range("L132").Value = 13950
$ws.Range("M132").Value = -9843.800000000001
$ws.Range("N132").Value = -19010
$ws.Range("H135").Value = 943.03125
$ws.Range("I135").Value = 882.04
$ws.Range("J135").Value = 1160.8572
$ws.Range("K135").Value = 7938.36
$ws.Range("L135").Value = 10447.7148
$ws.Range("M135").Value = -5403.36
$ws.Range("N135").Value = -15517.7148
$ws.Range("H137").Value = 1612.6394
$ws.Range("I137").Value = 911.8182
$ws.Range("J137").Value = 3426.5293
$ws.Range("K137").Value = 2735.4546
$ws.Range("L137").Value = 10279.5879
$ws.Range("M137").Value = -185.4546
$ws.Range("N137").Value = -15379.5879
$ws.Range("H138").Value = 1795.9452
$ws.Range("I138").Value = 987
$ws.Range("J138").Value = 2243.4468
$ws.Range("K138").Value = 2961
$ws.Range("L138").Value = 6730.340400000001
$ws.Range("M138").Value = 2179
$ws.Range("N138").Value = -17010.3404

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1286.2
$ws.Range("I45").Value = 911.5294
$ws.Range("J45").Value = 3409.3333
$ws.Range("K45").Value = 911.5294
$ws.Range("L45").Value = 3409.3333
$ws.Range("M45").Value = -534.5294
$ws.Range("N45").Value = -4163.3333
$ws.Range("H61").Value = 3050
$ws.Range("I61").Value = 2348.375
$ws.Range("J61").Value = 3985.5
$ws.Range("K61").Value = 2348.375
$ws.Range("L61").Value = 3985.5
$ws.Range("M61").Value = -2136.375
$ws.Range("N61").Value = -4409.5
$ws.Range("H110").Value = 3711.5356
$ws.Range("I110").Value = 2885.7222
$ws.Range("K110").Value = 2885.7222
$ws.Range("M110").Value = -840.7222000000002
$ws.Range("H122").Value = 1300.75
$ws.Range("I122").Value = 1156.4286
$ws.Range("J122").Value = 1502.8
$ws.Range("K122").Value = 3469.2858
$ws.Range("L122").Value = 4508.4
$ws.Range("M122").Value = -1019.2858
$ws.Range("N122").Value = -9408.4
$ws.Range("H132").Value = 1904.5652
$ws.Range("I132").Value = 1274.6154
$ws.Range("J132").Value = 2723.5
$ws.Range("K132").Value = 3823.8462
$ws.Range("L132").Value = 8170.5
$ws.Range("M132").Value = -1293.8462
$ws.Range("N132").Value = -13230.5
$ws.Range("H136").Value = 3050
$ws.Range("I136").Value = 2348.375
$ws.Range("J136").Value = 3985.5
$ws.Range("K136").Value = 7045.125
$ws.Range("L136").Value = 11956.5
$ws.Range("M136").Value = -4495.125
$ws.Range("N136").Value = -17056.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 495.6129
$ws.Range("I94").Value = 435.3913
$ws.Range("J94").Value = 668.75
$ws.Range("K94").Value = 435.3913
$ws.Range("L94").Value = 668.75
$ws.Range("M94").Value = 15.6087
$ws.Range("N94").Value = -1570.75
$ws.Range("H107").Value = 2121.4565
$ws.Range("I107").Value = 2096.5642
$ws.Range("J107").Value = 2260.1428
$ws.Range("K107").Value = 2096.5642
$ws.Range("L107").Value = 2260.1428
$ws.Range("M107").Value = -176.5641999999998
$ws.Range("N107").Value = -6100.1428
$ws.Range("H134").Value = 1663.7188
$ws.Range("I134").Value = 1351.1666
$ws.Range("J134").Value = 2601.375
$ws.Range("K134").Value = 4053.4998
$ws.Range("L134").Value = 7804.125
$ws.Range("M134").Value = -1518.4998
$ws.Range("N134").Value = -12874.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3132.75
$ws.Range("I31").Value = 2224.3635
$ws.Range("J31").Value = 5131.2
$ws.Range("K31").Value = 2224.3635
$ws.Range("L31").Value = 5131.2
$ws.Range("M31").Value = -1929.3635
$ws.Range("N31").Value = -5721.2
$ws.Range("H34").Value = 3132.75
$ws.Range("I34").Value = 2224.3635
$ws.Range("J34").Value = 5131.2
$ws.Range("K34").Value = 2224.3635
$ws.Range("L34").Value = 5131.2
$ws.Range("M34").Value = -2022.3635
$ws.Range("N34").Value = -5535.2
$ws.Range("H132").Value = 3185.2942
$ws.Range("I132").Value = 2312.3333
$ws.Range("J132").Value = 4167.375
$ws.Range("K132").Value = 6936.999899999999
$ws.Range("L132").Value = 12502.125
$ws.Range("M132").Value = -4406.999899999999
$ws.Range("N132").Value = -17562.125
$ws.Range("H134").Value = 3799.2632
$ws.Range("I134").Value = 1893.6
$ws.Range("J134").Value = 5916.6665
$ws.Range("K134").Value = 5680.799999999999
$ws.Range("L134").Value = 17749.9995
$ws.Range("M134").Value = -3145.799999999999
$ws.Range("N134").Value = -22819.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 2300
$ws.Range("J48").Value = 2300
$ws.Range("L48").Value = 6900
$ws.Range("N48").Value = -7400
$ws.Range("H50").Value = 67
$ws.Range("I50").Value = 45
$ws.Range("K50").Value = 135
$ws.Range("M50").Value = 346
$ws.Range("H53").Value = 67
$ws.Range("I53").Value = 45
$ws.Range("K53").Value = 135
$ws.Range("M53").Value = 346
$ws.Range("H131").Value = 848.94446
$ws.Range("J131").Value = 891.4516
$ws.Range("L131").Value = 2674.3548
$ws.Range("N131").Value = -12754.3548

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3226965.5
$ws.Range("I122").Value = 5556773.5
$ws.Range("J122").Value = 1077.8462
$ws.Range("K122").Value = 16670320.5
$ws.Range("L122").Value = 3233.5386
$ws.Range("M122").Value = -16667870.5
$ws.Range("N122").Value = -8133.5386
$ws.Range("H132").Value = 2734.697
$ws.Range("I132").Value = 2030.55
$ws.Range("J132").Value = 3818
$ws.Range("K132").Value = 6091.65
$ws.Range("L132").Value = 11454
$ws.Range("M132").Value = -3561.65
$ws.Range("N132").Value = -16514

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3763.5386
$ws.Range("I40").Value = 2997.6
$ws.Range("J40").Value = 6316.6665
$ws.Range("K40").Value = 2997.6
$ws.Range("L40").Value = 6316.6665
$ws.Range("M40").Value = -2861.6
$ws.Range("N40").Value = -6588.6665
$ws.Range("H108").Value = 23054
$ws.Range("J108").Value = 23054
$ws.Range("L108").Value = 23054
$ws.Range("N108").Value = -30734
$ws.Range("H122").Value = 3512.1206
$ws.Range("I122").Value = 3350.4614
$ws.Range("J122").Value = 3843.9473
$ws.Range("K122").Value = 10051.3842
$ws.Range("L122").Value = 11531.8419
$ws.Range("M122").Value = -7601.3842
$ws.Range("N122").Value = -16431.8419
$ws.Range("H140").Value = 60429
$ws.Range("J140").Value = 60429
$ws.Range("L140").Value = 60429
$ws.Range("N140").Value = -70789

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 448.57144
$ws.Range("I107").Value = 306.66666
$ws.Range("J107").Value = 555
$ws.Range("K107").Value = 919.9999799999999
$ws.Range("L107").Value = 1665
$ws.Range("M107").Value = 1000.00002
$ws.Range("N107").Value = -5505
$ws.Range("H122").Value = 1581.4546
$ws.Range("I122").Value = 1436.2222
$ws.Range("J122").Value = 2235
$ws.Range("K122").Value = 4308.6666
$ws.Range("L122").Value = 6705
$ws.Range("M122").Value = -1858.6666
$ws.Range("N122").Value = -11605
$ws.Range("H132").Value = 2004.2222
$ws.Range("I132").Value = 983.12
$ws.Range("K132").Value = 2949.36
$ws.Range("M132").Value = -419.3600000000001
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
